# "Generate Report for Handback"
#
# The localization-status report is regenerated by CI whenever a file is
# handed back from translation. The file 0ff44a1d-9ba4-4a34-a697-4e186131fe52.md
# has now come back from both zh-cn and de-de in sync with en-US, so its
# status flips from "Ready for handoff" to "Handed back: in sync with en-US"
# everywhere it is reported (Overview + each language sheet), and the
# per-language detail sheets gain their "Latest Target File" / "Latest
# Handback File" links plus a real "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Handed back: in sync with en-US"
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0dfa79caea41a2dc5d65047e5bc9b3caa2b4626e/e2e/0ff44a1d-9ba4-4a34-a697-4e186131fe52.md", "", "", "0ff44a1d-9ba4-4a34-a697-4e186131fe52.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a9251f3ea3e49ace1a1579137d0bf7a597f7b481/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/4616bf90-a0bb-4ee6-b432-e92cbf088aab.a9251f3ea3e49ace1a1579137d0bf7a597f7b481.zh-cn.xlf", "", "", "0ff44a1d-9ba4-4a34-a697-4e186131fe52.594d167c415c2d190e50774abdbe280b3156684a.zh-cn.xlf") | Out-Null
$zhcn.Range("G2").Value = "2016-03-10 14:30:15"
$zhcn.Range("H2").Value = "Include"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Handed back: in sync with en-US"
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0dfa79caea41a2dc5d65047e5bc9b3caa2b4626e/e2e/0ff44a1d-9ba4-4a34-a697-4e186131fe52.md", "", "", "0ff44a1d-9ba4-4a34-a697-4e186131fe52.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e3c0712d6154d3da8c431078bc941020631a4fec/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/4616bf90-a0bb-4ee6-b432-e92cbf088aab.a9251f3ea3e49ace1a1579137d0bf7a597f7b481.de-de.xlf", "", "", "0ff44a1d-9ba4-4a34-a697-4e186131fe52.594d167c415c2d190e50774abdbe280b3156684a.de-de.xlf") | Out-Null
$dede.Range("G2").Value = "2016-03-10 14:30:22"
$dede.Range("H2").Value = "Include"
